# "add api import from excel"
# The single-sheet "document import" template goes from 4 columns
# (Tieu de, Tac gia, Nha xuat ban, Nam xuat ban) to a 10-column layout that
# keeps the existing fields but adds the extra metadata needed by the new
# import API, re-ordered so the table reads:
#   Ten tai lieu, Tac gia, Nha xuat ban, Nam xuat ban, So hieu phan loai,
#   The loai, Vi tri tren ke, Loai tai lieu, Lien ket truy cap, Trang thai

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has a single Excel Table ("Table3") bound to A1:D1048575.
# Grow it to A1:J1048575 (10 columns) - this also appends the extra
# tableColumn entries the table needs.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J1048575"))

# Re-write the header row (row 1) in the final column order. Because the
# header row drives the table's column names, this both relabels the
# existing columns and names the newly-inserted ones.
$ws.Range("A1").Value = "Tên tài liệu"
$ws.Range("B1").Value = "Tác giả"
$ws.Range("C1").Value = "Nhà xuất bản"
$ws.Range("D1").Value = "Năm xuất bản"
$ws.Range("E1").Value = "Số hiệu phân loại"
$ws.Range("F1").Value = "Thể loại"
$ws.Range("G1").Value = "Vị trí trên kệ"
$ws.Range("H1").Value = "Loại tài liệu"
$ws.Range("I1").Value = "Liên kết truy cập"
$ws.Range("J1").Value = "Trạng thái"

# Match the new column widths: columns A-G share one width, H/I keep the
# old column B/C widths, and J gets its own width.
$ws.Range("A1:G1").ColumnWidth = 23.333333333333336
$ws.Range("H1").ColumnWidth = 14
$ws.Range("I1").ColumnWidth = 20.333333333333336
$ws.Range("J1").ColumnWidth = 17.833333333333336

# Move the active selection the way the author left it.
$ws.Range("F6").Select()
